# hotfix double-encoding-issue by using triple-braces
#
# The template uses Handlebars-style placeholders such as {{AccountName__c}}.
# Downstream processing was double-encoding these, so every double-brace
# placeholder ("{{Name}}" / "}}") is widened to a triple-brace placeholder
# ("{{{Name}}}" / "}}}") across the whole sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("{{", "{{{")
$ws.Cells.Replace("}}", "}}}")
